# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff" and its timestamps were refreshed. Reflect that on
# all three sheets (Overview, zh-cn, de-de) and widen the status columns
# that now hold the longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff timestamps ---
$overview.Range("G2").Value = "2016-08-30 09:22:57"
$zhcn.Range("H2").Value = "2016-08-30 09:22:44"
$dede.Range("H2").Value = "2016-08-30 09:22:57"

# --- Widen the status columns to fit "Ready for handoff" ---
$overview.Columns(5).ColumnWidth = 16.3
$overview.Columns(6).ColumnWidth = 16.3
$zhcn.Columns(3).ColumnWidth = 16.3
$dede.Columns(3).ColumnWidth = 16.3
